# The deck's Design/theme was switched from the custom "Integral" theme to
# the stock PowerPoint "Office Theme" palette (Design > Themes > Office in
# the UI). That swap changes the twelve theme colour slots that live in the
# slide master's theme part (ppt/theme/theme1.xml) from the Integral values
# to the standard Office values, while leaving the font scheme and format
# scheme (fills/lines/effects) untouched, since both themes already share
# identical font/format schemes.
#
# PowerPoint's object model doesn't expose a "set theme XML" call; the
# supported automation surface for recolouring a theme is to set each slot
# on ThemeColorScheme individually (Dark1, Light1, Dark2, Light2, Accent1-6,
# Hyperlink, FollowedHyperlink -> indices 1-12).

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

function ToRGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme colour scheme (standard PowerPoint default palette).
$cs.Item(1).RGB  = ToRGB 0x00 0x00 0x00   # Dark 1    (dk1)      000000
$cs.Item(2).RGB  = ToRGB 0xFF 0xFF 0xFF   # Light 1   (lt1)      FFFFFF
$cs.Item(3).RGB  = ToRGB 0x44 0x54 0x6A   # Dark 2    (dk2)      44546A
$cs.Item(4).RGB  = ToRGB 0xE7 0xE6 0xE6   # Light 2   (lt2)      E7E6E6
$cs.Item(5).RGB  = ToRGB 0x5B 0x9B 0xD5   # Accent 1             5B9BD5
$cs.Item(6).RGB  = ToRGB 0xED 0x7D 0x31   # Accent 2             ED7D31
$cs.Item(7).RGB  = ToRGB 0xA5 0xA5 0xA5   # Accent 3             A5A5A5
$cs.Item(8).RGB  = ToRGB 0xFF 0xC0 0x00   # Accent 4             FFC000
$cs.Item(9).RGB  = ToRGB 0x44 0x72 0xC4   # Accent 5             4472C4
$cs.Item(10).RGB = ToRGB 0x70 0xAD 0x47   # Accent 6             70AD47
$cs.Item(11).RGB = ToRGB 0x05 0x63 0xC1   # Hyperlink            0563C1
$cs.Item(12).RGB = ToRGB 0x95 0x4F 0x72   # Followed Hyperlink   954F72
